$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume 1h) values are numeric-looking text;
# force text format so Excel does not convert them to real numbers,
# matching the original inline-string text storage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "275.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.68%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.01%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.837"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06396"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.919"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.22%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.203"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.16%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8781"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1518"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05031"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.19%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07577"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.65%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02965"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.13%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08996"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001568"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.60%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006407"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006186"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.36%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.28%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.36%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1341"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.914"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.21%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001177"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.06%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003856"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-12.46%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.06%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "13.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04125"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.42%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006812"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.92%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.66%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002193"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.25%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01147"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005200"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.23%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.02003"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.02%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.650"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-47.11%"
